$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap match data (columns F:V) between rows 159 and 160 ---
$ws.Range("F159").Value = "Ankaragucu"
$ws.Range("G159").Value = 0
$ws.Range("H159").Value = "Hatayspor"
$ws.Range("I159").Value = 0
$ws.Range("J159").Value = 1.97
$ws.Range("K159").Value = "14/12/2023 18:13"
$ws.Range("L159").Value = 2.27
$ws.Range("M159").Value = "21/12/2023 17:56"
$ws.Range("N159").Value = 3.67
$ws.Range("O159").Value = "14/12/2023 18:13"
$ws.Range("P159").Value = 3.45
$ws.Range("Q159").Value = "21/12/2023 17:59"
$ws.Range("R159").Value = 3.82
$ws.Range("S159").Value = "14/12/2023 18:13"
$ws.Range("T159").Value = 3.36
$ws.Range("U159").Value = "21/12/2023 17:59"
$ws.Range("V159").Value = "https://www.betexplorer.com/football/turkey/super-lig/ankaragucu-hatayspor/Iw6iNEIm/"

$ws.Range("F160").Value = "Besiktas"
$ws.Range("G160").Value = 1
$ws.Range("H160").Value = "Alanyaspor"
$ws.Range("I160").Value = 3
$ws.Range("J160").Value = 1.43
$ws.Range("K160").Value = "14/12/2023 18:12"
$ws.Range("L160").Value = 1.64
$ws.Range("M160").Value = "21/12/2023 17:55"
$ws.Range("N160").Value = 4.9
$ws.Range("O160").Value = "14/12/2023 18:12"
$ws.Range("P160").Value = 4.32
$ws.Range("Q160").Value = "21/12/2023 17:55"
$ws.Range("R160").Value = 7
$ws.Range("S160").Value = "14/12/2023 18:12"
$ws.Range("T160").Value = 5.26
$ws.Range("U160").Value = "21/12/2023 17:57"
$ws.Range("V160").Value = "https://www.betexplorer.com/football/turkey/super-lig/besiktas-alanyaspor/l41mOf3s/"

# --- Swap match data (columns F:V) between rows 167 and 168 ---
$ws.Range("F167").Value = "Alanyaspor"
$ws.Range("G167").Value = 3
$ws.Range("H167").Value = "Samsunspor"
$ws.Range("I167").Value = 1
$ws.Range("J167").Value = 2.59
$ws.Range("K167").Value = "21/12/2023 18:12"
$ws.Range("L167").Value = 2.74
$ws.Range("M167").Value = "25/12/2023 14:53"
$ws.Range("N167").Value = 3.46
$ws.Range("O167").Value = "21/12/2023 18:12"
$ws.Range("P167").Value = 3.24
$ws.Range("Q167").Value = "25/12/2023 14:55"
$ws.Range("R167").Value = 2.76
$ws.Range("S167").Value = "21/12/2023 18:12"
$ws.Range("T167").Value = 2.83
$ws.Range("U167").Value = "25/12/2023 14:51"
$ws.Range("V167").Value = "https://www.betexplorer.com/football/turkey/super-lig/alanyaspor-samsunspor/Q9dCIjXI/"

$ws.Range("F168").Value = "Kasimpasa"
$ws.Range("G168").Value = 2
$ws.Range("H168").Value = "Rizespor"
$ws.Range("I168").Value = 2
$ws.Range("J168").Value = 1.93
$ws.Range("K168").Value = "21/12/2023 15:12"
$ws.Range("L168").Value = 2.45
$ws.Range("M168").Value = "25/12/2023 14:58"
$ws.Range("N168").Value = 3.8
$ws.Range("O168").Value = "21/12/2023 15:12"
$ws.Range("P168").Value = 3.68
$ws.Range("Q168").Value = "25/12/2023 14:58"
$ws.Range("R168").Value = 3.87
$ws.Range("S168").Value = "21/12/2023 15:12"
$ws.Range("T168").Value = 2.87
$ws.Range("U168").Value = "25/12/2023 14:58"
$ws.Range("V168").Value = "https://www.betexplorer.com/football/turkey/super-lig/kasimpasa-rizespor/nJc8JWHC/"

# --- Append new row 171 (new match row), copying formatting from the last data row ---
$ws.Range("A170:V170").Copy($ws.Range("A171:V171"))

$ws.Range("A171").Value = 170
$ws.Range("B171").Value = "turkey"
$ws.Range("C171").Value = "super-lig"
$ws.Range("D171").Value = "2023-2024"
$ws.Range("E171").Value = 45296.625
$ws.Range("F171").Value = "Gaziantep"
$ws.Range("G171").Value = 2
$ws.Range("H171").Value = "Pendikspor"
$ws.Range("I171").Value = 2
$ws.Range("J171").Value = 2.04
$ws.Range("K171").Value = "28/12/2024 15:42"
$ws.Range("L171").Value = 1.95
$ws.Range("M171").Value = "05/01/2024 14:56"
$ws.Range("N171").Value = 3.64
$ws.Range("O171").Value = "28/12/2024 15:42"
$ws.Range("P171").Value = 3.67
$ws.Range("Q171").Value = "05/01/2024 14:56"
$ws.Range("R171").Value = 3.62
$ws.Range("S171").Value = "28/12/2024 15:42"
$ws.Range("T171").Value = 4.1
$ws.Range("U171").Value = "05/01/2024 14:56"
$ws.Range("V171").Value = "https://www.betexplorer.com/football/turkey/super-lig/gaziantep-pendikspor/bTuUgDgP/"
